$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 24.40000000000038
$ws.Range("H2").Value = [double]"1.162537198560373e-16"
$ws.Range("I2").Value = 0.000233071194531087
$ws.Range("K2").Value = 57.81930919823225
$ws.Range("L2").Value = "[52.99210135602885, 62.64651704043566]"
$ws.Range("O2").Value = 1.62897396852804
$ws.Range("P2").Value = "[1.5409213215805782, 1.7170266154755023]"
$ws.Range("S2").Value = 53.76112404611962
$ws.Range("T2").Value = "[50.56567228127581, 56.95657581096343]"
$ws.Range("W2").Value = 18.07407407407435
$ws.Range("X2").Value = 17.7321321321324
$ws.Range("Y2").Value = 18.4160160160163

# Row 3 updates
$ws.Range("E3").Value = 23.84000000000029
$ws.Range("H3").Value = [double]"1.162537198560373e-16"
$ws.Range("K3").Value = 56.41555854058123
$ws.Range("L3").Value = "[49.47358686502699, 63.35753021613547]"
$ws.Range("O3").Value = 0.5723422051585008
$ws.Range("P3").Value = "[0.45913165908319353, 0.685552751233808]"
$ws.Range("S3").Value = 54.03091523078979
$ws.Range("T3").Value = "[50.4505672488077, 57.61126321277189]"
$ws.Range("W3").Value = 21.66838838838865
$ws.Range("X3").Value = 21.2388388388391
$ws.Range("Y3").Value = 22.0979379379382
